$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Myoc"
$ws.Cells.Item(2, 3).Value = "Fzd7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.349434
$ws.Cells.Item(2, 8).Value = 1.048302
$ws.Cells.Item(2, 9).Value = 0.0009963999680650763
$ws.Cells.Item(2, 10).Value = 0.0009963999680650763
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.553279333333334
$ws.Cells.Item(2, 14).Value = 7.659838000000001
$ws.Cells.Item(2, 15).Value = 0.1645043904057808
$ws.Cells.Item(2, 16).Value = 0.1645043904057808
$ws.Cells.Item(2, 17).Value = 0.8922026105640002
$ws.Cells.Item(2, 18).Value = 8.029823495076002
$ws.Cells.Item(2, 19).Value = 0.0001639121693468848
$ws.Cells.Item(2, 20).Value = 0.0001639121693468848

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Myoc"
$ws.Cells.Item(3, 3).Value = "Fzd7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.349434
$ws.Cells.Item(3, 8).Value = 1.048302
$ws.Cells.Item(3, 9).Value = 0.0009963999680650763
$ws.Cells.Item(3, 10).Value = 0.0009963999680650763
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.058662
$ws.Cells.Item(3, 14).Value = 24.175986
$ws.Cells.Item(3, 15).Value = 0.5192088709172035
$ws.Cells.Item(3, 16).Value = 0.5192088709172035
$ws.Cells.Item(3, 17).Value = 2.815970497308
$ws.Cells.Item(3, 18).Value = 25.343734475772
$ws.Cells.Item(3, 19).Value = 0.0005173397024010059
$ws.Cells.Item(3, 20).Value = 0.0005173397024010059

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Myoc"
$ws.Cells.Item(4, 3).Value = "Fzd7"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.349434
$ws.Cells.Item(4, 8).Value = 1.048302
$ws.Cells.Item(4, 9).Value = 0.0009963999680650763
$ws.Cells.Item(4, 10).Value = 0.0009963999680650763
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.909099333333334
$ws.Cells.Item(4, 14).Value = 14.727298
$ws.Cells.Item(4, 15).Value = 0.3162867386770157
$ws.Cells.Item(4, 16).Value = 0.3162867386770157
$ws.Cells.Item(4, 17).Value = 1.715406216444
$ws.Cells.Item(4, 18).Value = 15.438655947996
$ws.Cells.Item(4, 19).Value = 0.0003151480963171855
$ws.Cells.Item(4, 20).Value = 0.0003151480963171855

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Myoc"
$ws.Cells.Item(5, 3).Value = "Fzd7"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 348.977468
$ws.Cells.Item(5, 8).Value = 1046.932404
$ws.Cells.Item(5, 9).Value = 0.9950981815468188
$ws.Cells.Item(5, 10).Value = 0.9950981815468188
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.553279333333334
$ws.Cells.Item(5, 14).Value = 7.659838000000001
$ws.Cells.Item(5, 15).Value = 0.1645043904057808
$ws.Cells.Item(5, 16).Value = 0.1645043904057808
$ws.Cells.Item(5, 17).Value = 891.0369568433947
$ws.Cells.Item(5, 18).Value = 8019.332611590552
$ws.Cells.Item(5, 19).Value = 0.1636980197492604
$ws.Cells.Item(5, 20).Value = 0.1636980197492604

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Myoc"
$ws.Cells.Item(6, 3).Value = "Fzd7"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 348.977468
$ws.Cells.Item(6, 8).Value = 1046.932404
$ws.Cells.Item(6, 9).Value = 0.9950981815468188
$ws.Cells.Item(6, 10).Value = 0.9950981815468188
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 8.058662
$ws.Cells.Item(6, 14).Value = 24.175986
$ws.Cells.Item(6, 15).Value = 0.5192088709172035
$ws.Cells.Item(6, 16).Value = 0.5192088709172035
$ws.Cells.Item(6, 17).Value = 2812.291460227816
$ws.Cells.Item(6, 18).Value = 25310.62314205034
$ws.Cells.Item(6, 19).Value = 0.5166638032926861
$ws.Cells.Item(6, 20).Value = 0.5166638032926861

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Myoc"
$ws.Cells.Item(7, 3).Value = "Fzd7"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 348.977468
$ws.Cells.Item(7, 8).Value = 1046.932404
$ws.Cells.Item(7, 9).Value = 0.9950981815468188
$ws.Cells.Item(7, 10).Value = 0.9950981815468188
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 4.909099333333334
$ws.Cells.Item(7, 14).Value = 14.727298
$ws.Cells.Item(7, 15).Value = 0.3162867386770157
$ws.Cells.Item(7, 16).Value = 0.3162867386770157
$ws.Cells.Item(7, 17).Value = 1713.165055507155
$ws.Cells.Item(7, 18).Value = 15418.48549956439
$ws.Cells.Item(7, 19).Value = 0.3147363585048722
$ws.Cells.Item(7, 20).Value = 0.3147363585048722

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Myoc"
$ws.Cells.Item(8, 3).Value = "Fzd7"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 1.369616666666667
$ws.Cells.Item(8, 8).Value = 4.10885
$ws.Cells.Item(8, 9).Value = 0.00390541848511611
$ws.Cells.Item(8, 10).Value = 0.00390541848511611
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 2.553279333333334
$ws.Cells.Item(8, 14).Value = 7.659838000000001
$ws.Cells.Item(8, 15).Value = 0.1645043904057808
$ws.Cells.Item(8, 16).Value = 0.1645043904057808
$ws.Cells.Item(8, 17).Value = 3.497013929588889
$ws.Cells.Item(8, 18).Value = 31.4731253663
$ws.Cells.Item(8, 19).Value = 0.0006424584871734935
$ws.Cells.Item(8, 20).Value = 0.0006424584871734934

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Myoc"
$ws.Cells.Item(9, 3).Value = "Fzd7"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 1.369616666666667
$ws.Cells.Item(9, 8).Value = 4.10885
$ws.Cells.Item(9, 9).Value = 0.00390541848511611
$ws.Cells.Item(9, 10).Value = 0.00390541848511611
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 8.058662
$ws.Cells.Item(9, 14).Value = 24.175986
$ws.Cells.Item(9, 15).Value = 0.5192088709172035
$ws.Cells.Item(9, 16).Value = 0.5192088709172035
$ws.Cells.Item(9, 17).Value = 11.03727778623333
$ws.Cells.Item(9, 18).Value = 99.33550007610002
$ws.Cells.Item(9, 19).Value = 0.00202772792211631
$ws.Cells.Item(9, 20).Value = 0.002027727922116311

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Myoc"
$ws.Cells.Item(10, 3).Value = "Fzd7"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.369616666666667
$ws.Cells.Item(10, 8).Value = 4.10885
$ws.Cells.Item(10, 9).Value = 0.00390541848511611
$ws.Cells.Item(10, 10).Value = 0.00390541848511611
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 4.909099333333334
$ws.Cells.Item(10, 14).Value = 14.727298
$ws.Cells.Item(10, 15).Value = 0.3162867386770157
$ws.Cells.Item(10, 16).Value = 0.3162867386770157
$ws.Cells.Item(10, 17).Value = 6.723584265255556
$ws.Cells.Item(10, 18).Value = 60.51225838730001
$ws.Cells.Item(10, 19).Value = 0.001235232075826305
$ws.Cells.Item(10, 20).Value = 0.001235232075826305
